# "dsm matrix permutated styled"
# The DSM (Design Structure Matrix) body rows had a thin box border around
# every data cell (row index column A2:A16 plus the data columns B2:G16).
# After permutating/re-styling the matrix, that thin border is removed from
# the whole data body, leaving the cells borderless (the bold "thick"
# border on the header row, A1's center style, etc. are untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data body: row index column (A2:A16) + the six data columns (B2:G16).
# Both style groups currently reference the "thin" border; clearing the
# border on the whole block removes it for every cell that used it.
$body = $ws.Range("A2:G16")
$body.Borders.LineStyle = -4142   # xlLineStyleNone
